$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SKU numbers in column A to text values prefixed with "Hello"
$ids = @(2720010067, 2720010075, 2720010083, 2720010091, 2720010265, 2720010273, 2720010281, 2720010299, 2720010463)
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "Hello" + $ids[$i]
}

# Update the sheet view: scroll back to A1, select A10 instead of E1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A10").Select()
